$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Category" header cell in A1 and give it the same direct
# formatting (bold font, borders, centered alignment) as the rest of the
# header row (e.g. B1).
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

# The data cells A2:A46 previously shared the header's formatting; clear it
# so they match the (unstyled) formatting used by the rest of the data rows,
# e.g. B2, while keeping their existing text values untouched.
$ws.Range("B2").Copy()
$ws.Range("A2:A46").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
